$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.443.88'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.666.94'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +0.98%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.20'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +1.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3951'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +1.14%  '

$ws.Range("E8").Value = '  +0.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.06'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +5.87%  '

$ws.Range("E10").Value = '  +2.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.000'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  -0.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08584'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +0.86%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.40'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -0.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.326'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +1.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001341'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +3.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.890'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +4.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.664.87'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +1.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.53'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +0.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06964'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +0.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.55'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -3.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.007'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +0.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9991'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -0.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.73'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -1.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.426.18'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +0.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.434'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +3.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.010'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +8.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.52'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -0.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.79'
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '142.46'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -0.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.383'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +0.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.052'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -7.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.522'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +3.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.849.69'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +1.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.062'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +6.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08247'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +1.65%  '

$ws.Range("B36").Value = 'VeChain'

$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03003'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +1.29%  '

$ws.Range("B37").Value = 'FraxShare'

$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.15'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +10.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.786'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -3.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2767'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +1.77%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09258'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -0.79%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7737'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +0.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.82'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +4.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.448'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -2.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.62'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +2.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7128'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +2.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.533'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +0.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.144'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +0.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9993'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -0.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08453'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  -0.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.43'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +1.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.453'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +14.37%  '
